$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data for 2024-06-21 (Friday), matching the style of the
# preceding rows (date format on A/B, wrapped text on E).
$ws.Range("A40").Value = 45464
$ws.Range("A40").NumberFormat = "d-mmm"

$ws.Range("B40").Value = "F"
$ws.Range("B40").NumberFormat = "d-mmm"

$ws.Range("C40").Value = 0.5

$ws.Range("E40").Value = "Finished creating studentgrades_prof for test group"
$ws.Range("E40").WrapText = $true

# Update the active selection to follow the newly added row, as in the
# saved workbook.
$ws.Range("E40").Select() | Out-Null
